$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matching original inlineStr formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = '60.335.73'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '2.595.51'
$ws.Range("E3").Value = '  -2.71%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '509.37'
$ws.Range("E5").Value = '  -0.38%  '
$ws.Range("D6").Value = '153.85'
$ws.Range("E6").Value = '  -2.09%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  -2.61%  '
$ws.Range("D9").Value = '2.603.22'
$ws.Range("E9").Value = '  -2.26%  '
$ws.Range("D10").Value = '6.68'
$ws.Range("E10").Value = '  +4.83%  '
$ws.Range("D11").Value = '0.104'
$ws.Range("E11").Value = '  -1.16%  '
$ws.Range("E12").Value = '  -0.61%  '
$ws.Range("E13").Value = '  +1.66%  '
$ws.Range("D14").Value = '3.049.07'
$ws.Range("E14").Value = '  -2.29%  '
$ws.Range("D15").Value = '60.312.04'
$ws.Range("E15").Value = '  -1.05%  '
$ws.Range("D16").Value = '21.47'
$ws.Range("E16").Value = '  -1.45%  '
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = '2.598.55'
$ws.Range("E18").Value = '  -2.49%  '
$ws.Range("D19").Value = '4.74'
$ws.Range("E19").Value = '  -1.16%  '
$ws.Range("D20").Value = '353.32'
$ws.Range("E20").Value = '  +1.33%  '
$ws.Range("D21").Value = '10.50'
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("E22").Value = '  -0.70%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = '60.33'
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("D25").Value = '0.420'
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = '0.0₃0835'
$ws.Range("E28").Value = '  -2.77%  '
$ws.Range("D29").Value = '7.34'
$ws.Range("E29").Value = '  -2.49%  '
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").Value = '19.33'
$ws.Range("E31").Value = '  -0.80%  '
$ws.Range("D32").Value = '151.48'
$ws.Range("E32").Value = '  -4.25%  '
$ws.Range("E33").Value = '  -0.95%  '
$ws.Range("D34").Value = '5.72'
$ws.Range("E34").Value = '  +0.34%  '
$ws.Range("D35").Value = '3.99'
$ws.Range("E35").Value = '  -1.17%  '
$ws.Range("E36").Value = '  -2.96%  '
$ws.Range("E37").Value = '  +4.28%  '
$ws.Range("E38").Value = '  -2.93%  '
$ws.Range("D39").Value = '36.21'
$ws.Range("E39").Value = '  +2.43%  '
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").Value = '0.839'
$ws.Range("E40").Value = '  -2.44%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '3.75'
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").Value = '294.54'
$ws.Range("E42").Value = '  -4.83%  '
$ws.Range("E43").Value = '  -0.46%  '
$ws.Range("E44").Value = '  -4.04%  '
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("D47").Value = '19.62'
$ws.Range("E47").Value = '  -1.37%  '
$ws.Range("D48").Value = '4.88'
$ws.Range("E48").Value = '  -0.38%  '
$ws.Range("E49").Value = '  -1.26%  '
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = '1.988.36'
$ws.Range("E51").Value = '  -2.54%  '
